$d = $word.ActiveDocument

# 1. Generated-on sentence: date, time, and ARG version "develop" -> "1.1.7-RC4"
#    (do this first, since it also contains the "2021-08-26" substring that the
#    title-block replacement below would otherwise also catch)
$d.Content.Find.Execute("This document was generated on 2021-08-26, 08:26:23 with the Automatic Report Generator (ARG) version ""develop"" on the Linux system runner-ed2dce3a-project-18732201-concurrent-0.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "This document was generated on 2021-08-31, 15:49:27 with the Automatic Report Generator (ARG) version ""1.1.7-RC4"" on the Linux system runner-ed2dce3a-project-18732201-concurrent-0.",
                         2)

# 2. Title block date: 2021-08-26 -> 2021-08-31 (the standalone date after "root")
$d.Content.Find.Execute("2021-08-26", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2021-08-31", 2)

# 3. Remove numbering prefixes from headings
$d.Content.Find.Execute("1. Table of Contents", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Table of Contents", 2)

$d.Content.Find.Execute("2. List of Figures", $true, $false, $false, $false, $false,
                         $true, 1, $false, "List of Figures", 2)

$d.Content.Find.Execute("3. List of Tables", $true, $false, $false, $false, $false,
                         $true, 1, $false, "List of Tables", 2)

$d.Content.Find.Execute("4. Hello, world! Chapter!", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hello, world! Chapter!", 2)

$d.Content.Find.Execute("4.1. Hello, world! Section!", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hello, world! Section!", 2)

$d.Content.Find.Execute("4.1.1. Hello, world! Subsection!", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hello, world! Subsection!", 2)

$d.Content.Find.Execute("4.1.1.1. Hello, world! Subsubsection!", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hello, world! Subsubsection!", 2)

$d.Content.Find.Execute("4.1.2. Evidence description", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Evidence description", 2)

# Add a "_GoBack" bookmark right after the "Evidence description" run (zero-length,
# at the end of that paragraph's run, matching the expected output).
$r = $d.Content
$r.Find.Execute("Evidence description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
